# 26-Apr-2021, end of day update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update opening balance carried forward
$ws.Range("E2").Value = 685500

# Update row 3 entry: new date (26-Apr-2021) and adjusted Wages Expense amount
$ws.Range("A3").Value = 44312
$ws.Range("D3").Formula = "=60000"

# Update row 4 (TRANSFER BCA) expense total for the day
$ws.Range("D4").Formula = "=2877500+1537000+3649500"

# Remove all the now-obsolete daily transaction entries (rows 5 through 42),
# leaving just the running balance (column E) formulas in place.
$ws.Range("A5:D42").Clear()

# Reset view to top of the sheet and select D4, matching the new working cell.
[void]$ws.Range("D4").Select()
